$d = $word.ActiveDocument

# --- Change 1: "344+54=398" paragraph becomes ">>>" and a new
#     paragraph right after (currently empty) gets a run with "398",
#     carrying the same Arial run formatting as the original text. ---

# Find the paragraph that currently holds "344+54=398" and the
# (empty) paragraph right after it.
$count = $d.Paragraphs.Count
$srcParaIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $paraText = [string]$d.Paragraphs.Item($i).Range.Text
    if ($paraText -like "*344+54=398*") {
        $srcParaIndex = $i
        break
    }
}

if ($srcParaIndex -gt 0) {
    $srcPara = $d.Paragraphs.Item($srcParaIndex)
    $nextPara = $d.Paragraphs.Item($srcParaIndex + 1)

    # Copy a run from the source text so the pasted run picks up the
    # exact same direct run formatting (rFonts ascii/eastAsia/hAnsi/cs).
    $fmtRange = $d.Range($srcPara.Range.Start, $srcPara.Range.Start + 3)
    $fmtRange.Copy()
    $nextPara.Range.Paste()

    # Fix up the pasted text to be "398" instead of the copied "344".
    $pastedRange = $d.Range($nextPara.Range.Start, $nextPara.Range.Start + 3)
    $pastedRange.Text = "398"
}

# Now turn the original "344+54=398" text into ">>>"
$d.Content.Find.Execute("344+54=398", $true, $false, $false, $false, $false, `
                         $true, 1, $false, ">>>", 2)

# --- Change 2: merge "Serious exe" + "rcises" runs into a single
#     run reading "Serious exercises" (identical formatting already). ---
$d.Content.Find.Execute("Serious exercises", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Serious exercises", 2)

# --- Change 3: merge " ... stored unde" + "r the " runs into a single
#     run reading " from the list of items stored under the ". ---
$d.Content.Find.Execute(" from the list of items stored under the ", $true, `
                         $false, $false, $false, $false, $true, 1, $false, `
                         " from the list of items stored under the ", 2)
